# Update the workbook's build/version string from
#   "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# to
#   "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# across the "About" sheet (A2, A6) and the "build_version" column (S) of the
# "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Moranbah North Coal Mine, Australia, M0074, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds "build_version" for every data row (2 through the last used row).
$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$versionRange = $wsData.Range($wsData.Cells.Item(2, 19), $wsData.Cells.Item($lastRow, 19))
$versionRange.Value = $newVersion
